# automa.xlsx - "Add files via upload" edit
# Adds an "Area" column (MESA & BAR / COCINA) to the Tabla1 table on the
# "Personas" sheet, corrects four employee-code values in column A, and
# updates the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Personas")

# ---------------------------------------------------------------------
# 1) Fix a handful of employee codes in column A (rows 5, 14, 17, 23)
# ---------------------------------------------------------------------
$ws.Range("A5").Value = 2801
$ws.Range("A14").Value = 4863
$ws.Range("A17").Value = 5002
$ws.Range("A23").Value = 2007

# ---------------------------------------------------------------------
# 2) Grow the "Tabla1" table with a new "Area" column and populate it
# ---------------------------------------------------------------------
$lo = $ws.ListObjects.Item("Tabla1")
$newColumn = $lo.ListColumns.Add()

# Match the formatting already used by the other data columns (centered
# horizontally and vertically) by copying the format from column A.
$ws.Range("A2").Copy()
$ws.Range("C1:C31").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("C1").Value = "Area"

# Row 2 .. Row 31 -> Area, in sheet order (mirrors the source data)
$areas = @(
    "MESA & BAR",  # 2
    "MESA & BAR",  # 3
    "MESA & BAR",  # 4
    "MESA & BAR",  # 5
    "MESA & BAR",  # 6
    "MESA & BAR",  # 7
    "COCINA",      # 8
    "COCINA",      # 9
    "COCINA",      # 10
    "COCINA",      # 11
    "COCINA",      # 12
    "COCINA",      # 13
    "COCINA",      # 14
    "COCINA",      # 15
    "COCINA",      # 16
    "COCINA",      # 17
    "MESA & BAR",  # 18
    "MESA & BAR",  # 19
    "MESA & BAR",  # 20
    "MESA & BAR",  # 21
    "MESA & BAR",  # 22
    "MESA & BAR",  # 23
    "MESA & BAR",  # 24
    "MESA & BAR",  # 25
    "MESA & BAR",  # 26
    "MESA & BAR",  # 27
    "MESA & BAR",  # 28
    "MESA & BAR",  # 29
    "MESA & BAR",  # 30
    "MESA & BAR"   # 31
)

for ($i = 0; $i -lt $areas.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $areas[$i]
}

# ---------------------------------------------------------------------
# 3) Move the active selection, matching the saved workbook state
# ---------------------------------------------------------------------
$ws.Range("E5").Select() | Out-Null

Write-Host "Area column added and employee codes updated."
